$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E width change (16.3984375 -> 19.3984375)
# (ColumnWidth rounds to the nearest 1/6 char-width internally; 18.5 lands closest
# to the target stored width of 19.3984375)
$ws.Columns.Item(5).ColumnWidth = 18.5

# Selection change from E13 to A14
$ws.Range("A14").Select()

# E2 gets the shared string "different any occurence" (style stays the same fill)
$ws.Range("E2").Value = "different any occurence"

$greenColor = 5296274

# Row 8: C8:G8 fill becomes green (style index 2), values unchanged
$ws.Range("C8:G8").Interior.Color = $greenColor

# Row 9: C9, D9, F9 become green; E9, G9 stay red
$ws.Range("C9:D9").Interior.Color = $greenColor
$ws.Range("F9").Interior.Color = $greenColor

# Row 10: C10, F10 become green; D10, E10, G10 stay red
$ws.Range("C10").Interior.Color = $greenColor
$ws.Range("F10").Interior.Color = $greenColor

# E12: remove the shared string "different any occurence" -> empty, style stays red
$ws.Range("E12").Value = $null

# Row 13: C13, F13 become green; D13, E13, G13 stay red
$ws.Range("C13").Interior.Color = $greenColor
$ws.Range("F13").Interior.Color = $greenColor
